$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Maharashtra
$ws.Range("B3").Value = 12321
$ws.Range("C3").Value = 1498
$ws.Range("D3").Value = 420
$ws.Range("E3").Value = 10403

# Row 5 - Delhi
$ws.Range("B5").Value = 2916
$ws.Range("C5").Value = 295
$ws.Range("D5").Value = 187
$ws.Range("E5").Value = 2434

# Row 7 - Karnataka
$ws.Range("B7").Value = 1242
$ws.Range("C7").Value = 118
$ws.Range("D7").Value = 14
$ws.Range("E7").Value = 1110

# Row 8 - Rajasthan
$ws.Range("B8").Value = 1046
$ws.Range("E8").Value = 888

# Row 9 - Gujarat
$ws.Range("B9").Value = 938
$ws.Range("E9").Value = 821

# Row 10 - Telangana
$ws.Range("B10").Value = 735
$ws.Range("C10").Value = 57
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 667

# Row 11 - Haryana
$ws.Range("B11").Value = 766
$ws.Range("C11").Value = 64
$ws.Range("D11").Value = 33
$ws.Range("E11").Value = 669

# Row 12 - Tamil Nadu
$ws.Range("B12").Value = 650
$ws.Range("C12").Value = 118
$ws.Range("E12").Value = 514

# Row 13 - Andhra Pradesh
$ws.Range("B13").Value = 525
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 14
$ws.Range("E13").Value = 491

# Row 14 - Uttar Pradesh
$ws.Range("B14").Value = 387
$ws.Range("C14").Value = 218
$ws.Range("E14").Value = 167

# Row 15 - Madhya Pradesh
$ws.Range("B15").Value = 300
$ws.Range("C15").Value = 36
$ws.Range("E15").Value = 260

# Row 16 - Jammu and Kashmir
$ws.Range("B16").Value = 279
$ws.Range("C16").Value = 80
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = 187

# Row 17 - West Bengal
$ws.Range("B17").Value = 204
$ws.Range("C17").Value = 57
$ws.Range("E17").Value = 144

# Row 19 - Punjab
$ws.Range("B19").Value = 186
$ws.Range("E19").Value = 146

# Row 20 - Bihar
$ws.Range("B20").Value = 70
$ws.Range("E20").Value = 40

# Row 23 - Odisha
$ws.Range("C23").Value = 17
$ws.Range("E23").Value = 16

# Row 24 - Andaman and Nicobar Islands
$ws.Range("B24").Value = 35
$ws.Range("E24").Value = 20

# Row 25 - Chandigarh
$ws.Range("C25").Value = 2
$ws.Range("E25").Value = 29

# Row 26 - Ladakh
$ws.Range("B26").Value = 28
$ws.Range("E26").Value = 26

# Row 27 - Chhattisgarh
$ws.Range("C27").Value = 9
$ws.Range("E27").Value = 12

# Row 30 - Puducherry
$ws.Range("C30").Value = 6
$ws.Range("E30").Value = 1

# Row 33 - Dadra and Nagar Haveli
$ws.Range("C33").Value = 1
$ws.Range("E33").Value = 1

# Row 35 - Mizoram
$ws.Range("C35").Value = 1
$ws.Range("E35").Value = 0

# Row 38 - Sikkim
$ws.Range("B38").Value = 7
$ws.Range("E38").Value = 6

# Update selection on the sheet
$ws.Range("E1").Select()
